$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "0.999"
# or thousands-dotted numbers are not auto-converted to numeric types.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.316.70"
$ws.Range("E2").Value = "  -4.81%  "

$ws.Range("D3").Value = "3.256.39"
$ws.Range("E3").Value = "  -7.73%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "590.67"
$ws.Range("E5").Value = "  -5.34%  "

$ws.Range("D6").Value = "152.51"
$ws.Range("E6").Value = "  -12.38%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.250.32"
$ws.Range("E8").Value = "  -7.87%  "

$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  -10.57%  "

$ws.Range("E10").Value = "  -12.80%  "

$ws.Range("D11").Value = "6.86"
$ws.Range("E11").Value = "  -4.04%  "

$ws.Range("D12").Value = "0.508"
$ws.Range("E12").Value = "  -13.26%  "

$ws.Range("D13").Value = "38.57"
$ws.Range("E13").Value = "  -17.07%  "

$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -11.47%  "

$ws.Range("D15").Value = "3.773.90"
$ws.Range("E15").Value = "  -7.89%  "

$ws.Range("D16").Value = "67.349.13"
$ws.Range("E16").Value = "  -4.86%  "

$ws.Range("D17").Value = "548.41"
$ws.Range("E17").Value = "  -9.57%  "

$ws.Range("D18").Value = "3.253.34"
$ws.Range("E18").Value = "  -7.94%  "

$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  -13.60%  "

$ws.Range("E20").Value = "  -6.00%  "

$ws.Range("D21").Value = "15.23"
$ws.Range("E21").Value = "  -14.35%  "

$ws.Range("D22").Value = "0.769"
$ws.Range("E22").Value = "  -13.12%  "

$ws.Range("D23").Value = "7.96"
$ws.Range("E23").Value = "  -12.64%  "

$ws.Range("D24").Value = "85.77"
$ws.Range("E24").Value = "  -12.65%  "

$ws.Range("D25").Value = "13.60"
$ws.Range("E25").Value = "  -13.08%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -14.58%  "

$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  -10.04%  "

$ws.Range("D29").Value = "29.52"
$ws.Range("E29").Value = "  -12.68%  "

$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -17.04%  "

$ws.Range("E31").Value = "  -11.42%  "

$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -11.63%  "

$ws.Range("D33").Value = "547.88"
$ws.Range("E33").Value = "  -14.29%  "

$ws.Range("D34").Value = "6.66"
$ws.Range("E34").Value = "  -18.10%  "

$ws.Range("D35").Value = "5.78"
$ws.Range("E35").Value = "  -15.04%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").Value = "0.0450"
$ws.Range("E37").Value = "  -5.20%  "

$ws.Range("D38").Value = "53.69"
$ws.Range("E38").Value = "  -5.51%  "

$ws.Range("D39").Value = "0.0857"
$ws.Range("E39").Value = "  -14.14%  "

$ws.Range("D40").Value = "9.25"
$ws.Range("E40").Value = "  -14.58%  "

$ws.Range("E41").Value = "  -11.76%  "

$ws.Range("D42").Value = "2.934.07"
$ws.Range("E42").Value = "  -12.72%  "

$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -24.40%  "

$ws.Range("D44").Value = "0.263"
$ws.Range("E44").Value = "  -15.80%  "

$ws.Range("D45").Value = "0.0" + [string][char]8323 + "0588"
$ws.Range("E45").Value = "  -19.74%  "

$ws.Range("D46").Value = "26.62"
$ws.Range("E46").Value = "  -17.06%  "

$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  -20.17%  "

$ws.Range("D48").Value = "2.15"
$ws.Range("E48").Value = "  -15.77%  "

$ws.Range("D50").Value = "126.99"
$ws.Range("E50").Value = "  -4.40%  "

$ws.Range("E51").Value = "  -12.39%  "
